$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D12").Value = "2023년을 회고하며"
$ws.Range("E12").Value = "https://tensorflow.blog/2023/12/29/2023%eb%85%84%ec%9d%84-%ed%9a%8c%ea%b3%a0%ed%95%98%eb%a9%b0/"

$ws.Range("D51").Value = "[Linux] 리눅스 특정 폴더 찾기, find"
$ws.Range("E51").Value = "https://bskyvision.com/entry/Linux-%EB%A6%AC%EB%88%85%EC%8A%A4-%ED%8A%B9%EC%A0%95-%ED%8F%B4%EB%8D%94-%EC%B0%BE%EA%B8%B0-find"
